$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns remain text (matches original inlineStr cell type)
# so values like "1.001", "0.09865", "1.0000" are not reinterpreted as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.058.15'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.57%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.818.09'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.45%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.87'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.19%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4999'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -2.47%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3910'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.40%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.09865'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +26.15%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.48%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '40.88'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.02%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.438'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +3.08%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.79%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.002'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.19%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.818.90'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.05%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.284'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.57%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +5.51%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '92.36'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.61%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06643'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +1.75%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.17%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.21'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.79%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.958'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.45%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.111.73'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.47%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.27'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.246'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.75%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.22'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.89%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.029.37'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.75%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.69'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.82%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.413'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.76%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.96'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.20%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1064'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.37%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.034'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.84%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.571'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.40%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.620'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.10%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.06688'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -5.66%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02339'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.910'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +0.88%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2142'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.47%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.962'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -1.10%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '11.37'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.49%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6199'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.59%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.180'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.54%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.10%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.16'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.40%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5914'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.12%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.287'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.79%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.693'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.22%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.54'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.12%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.53%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.179'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.36%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06789'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.47%  '
